# Auto-generated COM-interop edit script.
# Adds the 'ProsthoWolf' backstory lines (rows 2-6) to 工作表1 / 工作表2
# and refreshes the selections accordingly.

$wb = $excel.ActiveWorkbook

# --- Sheet '工作表1' (TRANSLATE/DETECTLANGUAGE formula sheet) ---
$ws1 = $wb.Worksheets.Item("工作表1")

$ws1.Range("A2").Value = 'I refuse.
With your poor oral hygiene, your tooth would just decay until it fractures even if it is protected by a crown.'
$ws1.Range("A3").Value = 'Hey!
How dare you call my oral hygiene bad!'
$ws1.Range("A4").Value = 'Please do not let you food drop onto the floor.'
$ws1.Range("A5").Value = 'You place your dental implant too shallow!'
$ws1.Range("A6").Value = 'The emergence profile of your crown is ugly.'

$ws1.Rows.Item(2).RowHeight = 92.25
$ws1.Rows.Item(3).RowHeight = 46.15
$ws1.Rows.Item(4).RowHeight = 30.75
$ws1.Rows.Item(5).RowHeight = 30.75
$ws1.Rows.Item(6).RowHeight = 30.75

$ws1.Application.Goto($ws1.Range("A2:H6"))
$ws1.Range("A2:H6").Select()

# --- Sheet '工作表2' (plain literal translation grid) ---
$ws2 = $wb.Worksheets.Item("工作表2")

$ws2.Range("A2").Value = 'I refuse.
With your poor oral hygiene, your tooth would just decay until it fractures even if it is protected by a crown.'
$ws2.Range("B2").Value = 'en'
$ws2.Range("C2").Value = 'I refuse.
With your poor oral hygiene, your tooth would just decay until it fractures even if it is protected by a crown.'
$ws2.Range("D2").Value = '我拒絕。
口腔衛生這麼差，做了牙冠也是蛀掉。'
$ws2.Range("E2").Value = '私は拒否します。
口腔衛生状態が悪いと、歯冠で保護されていても、歯が骨折するまで虫歯になってしまいます。'
$ws2.Range("F2").Value = 'Me niego.
Con una mala higiene bucal, su diente simplemente se deterioraría hasta fracturarse, incluso si está protegido por una corona.'
$ws2.Range("G2").Value = 'Je refuse.
Avec votre mauvaise hygiène bucco-dentaire, votre dent ne ferait que se décomposer jusqu’à ce qu’elle se fracture, même si elle est protégée par une couronne.'
$ws2.Range("H2").Value = 'ฉันปฏิเสธ
ด้วยสุขอนามัยในช่องปากที่ไม่ดีฟันของคุณจะผุจนหักแม้ว่าจะได้รับการปกป้องด้วยครอบฟันก็ตาม'

$ws2.Range("A3").Value = 'Hey!
How dare you call my oral hygiene bad!'
$ws2.Range("B3").Value = 'en'
$ws2.Range("C3").Value = 'Hey!
How dare you call my oral hygiene bad!'
$ws2.Range("D3").Value = '欸！
說誰的口腔衛生差呢！?'
$ws2.Range("E3").Value = 'ねえ！
私の口腔衛生を悪いと呼ぶなんて!'
$ws2.Range("F3").Value = '¡Eh!
¡Cómo te atreves a decir que mi higiene bucal es mala!'
$ws2.Range("G3").Value = 'Hé!
Comment osez-vous dire que mon hygiène bucco-dentaire est mauvaise !'
$ws2.Range("H3").Value = 'หวัดดี!
คุณกล้าเรียกสุขอนามัยในช่องปากของฉันว่าไม่ดีได้อย่างไร!'

$ws2.Range("A4").Value = 'Please do not let you food drop onto the floor.'
$ws2.Range("B4").Value = 'en'
$ws2.Range("C4").Value = 'Please do not let you food drop onto the floor.'
$ws2.Range("D4").Value = '請不要讓食物掉到診間地板。'
$ws2.Range("E4").Value = '食べ物を床に落とさないでください。'
$ws2.Range("F4").Value = 'Por favor, no deje que la comida caiga al suelo.'
$ws2.Range("G4").Value = 'Veuillez ne pas laisser tomber votre nourriture sur le sol.'
$ws2.Range("H4").Value = 'โปรดอย่าปล่อยให้อาหารหล่นลงบนพื้น'

$ws2.Range("A5").Value = 'You place your dental implant too shallow!'
$ws2.Range("B5").Value = 'en'
$ws2.Range("C5").Value = 'You place your dental implant too shallow!'
$ws2.Range("D5").Value = '你植的太淺了！'
$ws2.Range("E5").Value = '歯科インプラントを浅くしすぎます!'
$ws2.Range("F5").Value = '¡Colocas tu implante dental demasiado poco profundo!'
$ws2.Range("G5").Value = 'Vous placez votre implant dentaire trop peu profond !'
$ws2.Range("H5").Value = 'คุณวางรากฟันเทียมตื้นเกินไป!'

$ws2.Range("A6").Value = 'The emergence profile of your crown is ugly.'
$ws2.Range("B6").Value = 'en'
$ws2.Range("C6").Value = 'The emergence profile of your crown is ugly.'
$ws2.Range("D6").Value = '是你的 emergence profile 醜~'
$ws2.Range("E6").Value = 'あなたの王冠の出現プロファイルは醜いです。'
$ws2.Range("F6").Value = 'El perfil de emergencia de su corona es feo.'
$ws2.Range("G6").Value = 'Le profil d’émergence de votre couronne est moche.'
$ws2.Range("H6").Value = 'โปรไฟล์การเกิดขึ้นของมงกุฎของคุณน่าเกลียด'

$ws2.Rows.Item(2).RowHeight = 48.4

$ws2.Application.Goto($ws2.Range("D7"))
$ws2.Range("D7").Select()

